$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.557.84"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.025.22"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "264.35"
$ws.Range("E5").Value = "  +6.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.620"
$ws.Range("E6").Value = "  -1.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.52"
$ws.Range("E8").Value = "  -6.96%  "
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0782"
$ws.Range("E10").Value = "  -2.60%  "
$ws.Range("E11").Value = "  -1.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.47"
$ws.Range("E12").Value = "  -3.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.327.79"
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.812"
$ws.Range("E14").Value = "  -4.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.95"
$ws.Range("E15").Value = "  -8.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.29"
$ws.Range("E16").Value = "  -3.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.035.16"
$ws.Range("E17").Value = "  +1.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.430.98"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.09"
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0849"
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.20"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.27"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.72"
$ws.Range("E23").Value = "  +8.14%  "
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("E25").Value = "  -1.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.63"
$ws.Range("E26").Value = "  +0.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.05"
$ws.Range("E27").Value = "  -3.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.77"
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("E29").Value = "  -10.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.35"
$ws.Range("E30").Value = "  +1.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.120"
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0656"
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.66"
$ws.Range("E33").Value = "  -3.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.56"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.41"
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.35"
$ws.Range("E38").Value = "  +2.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.26"
$ws.Range("E39").Value = "  -4.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.06"
$ws.Range("E40").Value = "  +4.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.23"
$ws.Range("E41").Value = "  +3.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0949"
$ws.Range("E42").Value = "  -3.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0215"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.409.18"
$ws.Range("E44").Value = "  +2.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.29"
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.80"
$ws.Range("E46").Value = "  -4.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.04"
$ws.Range("E47").Value = "  -0.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.14"
$ws.Range("E48").Value = "  -1.59%  "
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.219.33"
$ws.Range("E50").Value = "  +1.30%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.97"
$ws.Range("E51").Value = "  -2.14%  "
